$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 314, pushing existing rows 314:429 down to 315:430
$ws.Rows.Item(314).Insert()

# Populate the newly inserted row 314 with the new record's data
$ws.Cells.Item(314, 1).Value = 8
$ws.Cells.Item(314, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(314, 3).Value = "Coquimbo"
$ws.Cells.Item(314, 4).Value = Get-Date -Year 2023 -Month 8 -Day 8 -Hour 0 -Minute 0 -Second 0
$ws.Cells.Item(314, 5).Value = 4
$ws.Cells.Item(314, 6).Value = 100112031
$ws.Cells.Item(314, 7).Value = "Poroto verde"
$ws.Cells.Item(314, 8).Value = "Magnum"
$ws.Cells.Item(314, 9).Value = "Primera"
$ws.Cells.Item(314, 10).Value = 400
$ws.Cells.Item(314, 11).Value = 32000
$ws.Cells.Item(314, 12).Value = 33000
$ws.Cells.Item(314, 13).Value = 32500
$ws.Cells.Item(314, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(314, 15).Value = "Perú"
$ws.Cells.Item(314, 16).Value = 1300
$ws.Cells.Item(314, 17).Value = 25
$ws.Cells.Item(314, 18).Value = "Hortaliza"
